$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Development, 16-Dec-2021, 119 executed, 118 pass, 1 fail
$ws.Range("A5").Value = 44546
$ws.Range("B5").Value = "Development"
$ws.Range("C5").Value = 119
$ws.Range("D5").Value = 118
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "After execution all test cases pass"
$ws.Range("G5").Value = "Test cases iitially fail because of page load affected by network"

# Row 6: Production, 15-Dec-2021, 134 executed, 131 pass, 3 fail
$ws.Range("A6").Value = 44545
$ws.Range("B6").Value = "Production"
$ws.Range("C6").Value = 134
$ws.Range("D6").Value = 131
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = "After execution all test cases pass"
$ws.Range("G6").Value = "Test cases iitially fail because of page load affected by network"

# Copy formatting from row 4 down to rows 5 and 6
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A6:G6").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("A6").NumberFormat = $ws.Range("A4").NumberFormat

$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight
$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(4).RowHeight

# Update selection / view to match where the user ended up after the run
$ws.Range("E6").Select()
